$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3829.9
$ws.Range("I64").Value = 3287.375
$ws.Range("K64").Value = 3287.375
$ws.Range("M64").Value = -3039.375
$ws.Range("H67").Value = 3829.9
$ws.Range("I67").Value = 3287.375
$ws.Range("K67").Value = 3287.375
$ws.Range("M67").Value = -2429.375
$ws.Range("H100").Value = 1975.7727
$ws.Range("I100").Value = 1940.1
$ws.Range("K100").Value = 1940.1
$ws.Range("M100").Value = -1399.1
$ws.Range("H137").Value = 4112.3076
$ws.Range("J137").Value = 4304.5713
$ws.Range("L137").Value = 12913.7139
$ws.Range("N137").Value = -18013.7139
$ws.Range("H141").Value = 7151.213
$ws.Range("I141").Value = 6545.8
$ws.Range("K141").Value = 19637.4
$ws.Range("M141").Value = -14457.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1269.98
$ws.Range("I32").Value = 1178.266
$ws.Range("J32").Value = 2706.8333
$ws.Range("K32").Value = 1178.266
$ws.Range("L32").Value = 2706.8333
$ws.Range("M32").Value = -891.2660000000001
$ws.Range("N32").Value = -3280.8333
$ws.Range("H61").Value = 4764.4224
$ws.Range("I61").Value = 4729.108
$ws.Range("J61").Value = 4927.75
$ws.Range("K61").Value = 4729.108
$ws.Range("L61").Value = 4927.75
$ws.Range("M61").Value = -4517.108
$ws.Range("N61").Value = -5351.75
$ws.Range("H74").Value = 2226
$ws.Range("I74").Value = 1601.5454
$ws.Range("J74").Value = 3599.8
$ws.Range("K74").Value = 1601.5454
$ws.Range("L74").Value = 3599.8
$ws.Range("M74").Value = -727.5454
$ws.Range("N74").Value = -5347.8
$ws.Range("H77").Value = 2226
$ws.Range("I77").Value = 1601.5454
$ws.Range("J77").Value = 3599.8
$ws.Range("K77").Value = 8007.727
$ws.Range("L77").Value = 17999
$ws.Range("M77").Value = -3639.727
$ws.Range("N77").Value = -26735
$ws.Range("H132").Value = 3063.125
$ws.Range("I132").Value = 3293.5715
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 9880.7145
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -7350.7145
$ws.Range("N132").Value = -9410
$ws.Range("H136").Value = 4764.4224
$ws.Range("I136").Value = 4729.108
$ws.Range("J136").Value = 4927.75
$ws.Range("K136").Value = 14187.324
$ws.Range("L136").Value = 14783.25
$ws.Range("M136").Value = -11637.324
$ws.Range("N136").Value = -19883.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2877.7856
$ws.Range("I86").Value = 1898
$ws.Range("J86").Value = 3041.0833
$ws.Range("K86").Value = 1898
$ws.Range("L86").Value = 3041.0833
$ws.Range("M86").Value = -775
$ws.Range("N86").Value = -5287.0833
$ws.Range("H89").Value = 2877.7856
$ws.Range("I89").Value = 1898
$ws.Range("J89").Value = 3041.0833
$ws.Range("K89").Value = 9490
$ws.Range("L89").Value = 15205.4165
$ws.Range("M89").Value = -3874
$ws.Range("N89").Value = -26437.4165

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2045.7826
$ws.Range("I31").Value = 1692.25
$ws.Range("K31").Value = 1692.25
$ws.Range("M31").Value = -1397.25
$ws.Range("H34").Value = 2045.7826
$ws.Range("I34").Value = 1692.25
$ws.Range("K34").Value = 1692.25
$ws.Range("M34").Value = -1490.25
$ws.Range("H99").Value = 9748.8125
$ws.Range("I99").Value = 5895
$ws.Range("J99").Value = 14703.714
$ws.Range("K99").Value = 5895
$ws.Range("L99").Value = 14703.714
$ws.Range("M99").Value = -4397
$ws.Range("N99").Value = -17699.714
$ws.Range("H126").Value = 9748.8125
$ws.Range("I126").Value = 5895
$ws.Range("J126").Value = 14703.714
$ws.Range("K126").Value = 17685
$ws.Range("L126").Value = 44111.142
$ws.Range("M126").Value = -15215
$ws.Range("N126").Value = -49051.142
$ws.Range("H132").Value = 7772.44
$ws.Range("I132").Value = 7943.524
$ws.Range("K132").Value = 23830.572
$ws.Range("M132").Value = -21300.572
$ws.Range("H134").Value = 7274.7915
$ws.Range("I134").Value = 6540.4116
$ws.Range("J134").Value = 9058.286
$ws.Range("K134").Value = 19621.2348
$ws.Range("L134").Value = 27174.858
$ws.Range("M134").Value = -17086.2348
$ws.Range("N134").Value = -32244.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 313.45456
$ws.Range("J5").Value = 550.75
$ws.Range("L5").Value = 1652.25
$ws.Range("N5").Value = -1876.25
$ws.Range("H68").Value = 3221.375
$ws.Range("J68").Value = 3688.7693
$ws.Range("L68").Value = 11066.3079
$ws.Range("N68").Value = -12688.3079
$ws.Range("H71").Value = 3221.375
$ws.Range("J71").Value = 3688.7693
$ws.Range("L71").Value = 33198.9237
$ws.Range("N71").Value = -41310.9237
$ws.Range("H135").Value = 313.45456
$ws.Range("J135").Value = 550.75
$ws.Range("L135").Value = 4956.75
$ws.Range("N135").Value = -10026.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13221.272
$ws.Range("I70").Value = 16157
$ws.Range("J70").Value = 5392.6665
$ws.Range("K70").Value = 16157
$ws.Range("L70").Value = 5392.6665
$ws.Range("M70").Value = -15887
$ws.Range("N70").Value = -5932.6665
$ws.Range("H73").Value = 13221.272
$ws.Range("I73").Value = 16157
$ws.Range("J73").Value = 5392.6665
$ws.Range("K73").Value = 16157
$ws.Range("L73").Value = 5392.6665
$ws.Range("M73").Value = -15221
$ws.Range("N73").Value = -7264.6665
$ws.Range("H126").Value = 4764.2
$ws.Range("I126").Value = 4732.7646
$ws.Range("J126").Value = 4787.4346
$ws.Range("K126").Value = 14198.2938
$ws.Range("L126").Value = 14362.3038
$ws.Range("M126").Value = -11728.2938
$ws.Range("N126").Value = -19302.3038
$ws.Range("H132").Value = 3155.0322
$ws.Range("I132").Value = 3130.2334
$ws.Range("K132").Value = 9390.700199999999
$ws.Range("M132").Value = -6860.700199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3795.1765
$ws.Range("I7").Value = 3310.5
$ws.Range("J7").Value = 4958.4
$ws.Range("K7").Value = 3310.5
$ws.Range("L7").Value = 4958.4
$ws.Range("M7").Value = -3198.5
$ws.Range("N7").Value = -5182.4
$ws.Range("H16").Value = 2887.5908
$ws.Range("I16").Value = 2961.25
$ws.Range("K16").Value = 2961.25
$ws.Range("M16").Value = -2791.25
$ws.Range("H22").Value = 4556.7
$ws.Range("I22").Value = 4208.75
$ws.Range("J22").Value = 5948.5
$ws.Range("K22").Value = 4208.75
$ws.Range("L22").Value = 5948.5
$ws.Range("M22").Value = -3913.75
$ws.Range("N22").Value = -6538.5
$ws.Range("H27").Value = 4556.7
$ws.Range("I27").Value = 4208.75
$ws.Range("J27").Value = 5948.5
$ws.Range("K27").Value = 4208.75
$ws.Range("L27").Value = 5948.5
$ws.Range("M27").Value = -4101.75
$ws.Range("N27").Value = -6162.5
$ws.Range("H40").Value = 5246.4
$ws.Range("I40").Value = 3638.5
$ws.Range("K40").Value = 3638.5
$ws.Range("M40").Value = -3502.5
$ws.Range("H68").Value = 1705.7778
$ws.Range("I68").Value = 1499.6666
$ws.Range("J68").Value = 1808.8334
$ws.Range("K68").Value = 1499.6666
$ws.Range("L68").Value = 1808.8334
$ws.Range("M68").Value = -750.6666
$ws.Range("N68").Value = -3306.8334
$ws.Range("H71").Value = 1705.7778
$ws.Range("I71").Value = 1499.6666
$ws.Range("J71").Value = 1808.8334
$ws.Range("K71").Value = 7498.333000000001
$ws.Range("L71").Value = 9044.166999999999
$ws.Range("M71").Value = -3754.333000000001
$ws.Range("N71").Value = -16532.167
$ws.Range("H75").Value = 60173
$ws.Range("J75").Value = 60173
$ws.Range("L75").Value = 60173
$ws.Range("N75").Value = -62045
$ws.Range("H78").Value = 60173
$ws.Range("J78").Value = 60173
$ws.Range("L78").Value = 180519
$ws.Range("N78").Value = -189879
$ws.Range("H126").Value = 3795.1765
$ws.Range("I126").Value = 3310.5
$ws.Range("J126").Value = 4958.4
$ws.Range("K126").Value = 9931.5
$ws.Range("L126").Value = 14875.2
$ws.Range("M126").Value = -7461.5
$ws.Range("N126").Value = -19815.2
$ws.Range("H136").Value = 9482984
$ws.Range("I136").Value = 12003593
$ws.Range("J136").Value = 30701
$ws.Range("K136").Value = 36010779
$ws.Range("L136").Value = 92103
$ws.Range("M136").Value = -36008229
$ws.Range("N136").Value = -97203

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1738.5
$ws.Range("I107").Value = 854.5
$ws.Range("K107").Value = 2563.5
$ws.Range("M107").Value = -643.5
$ws.Range("H128").Value = 136475
$ws.Range("J128").Value = 136475
$ws.Range("L128").Value = 136475
$ws.Range("N128").Value = -146435
$ws.Range("H132").Value = 6228.436
$ws.Range("I132").Value = 6122.4688
$ws.Range("K132").Value = 18367.4064
$ws.Range("M132").Value = -15837.4064
$ws.Range("H136").Value = 17713.105
$ws.Range("I136").Value = 25882.166
$ws.Range("J136").Value = 9188.869000000001
$ws.Range("K136").Value = 77646.49800000001
$ws.Range("L136").Value = 27566.607
$ws.Range("M136").Value = -75096.49800000001
$ws.Range("N136").Value = -32666.607
